$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 44453
$ws.Range("A14").NumberFormat = "d-mmm"
$ws.Range("B14").Value = "4 Hours"
$ws.Range("C14").Value = "Self Learning: Watched videos on try catches in sql and transactions in sql, also watched videos on Delegates in c# and Eventhandlers"

$ws.Range("A15").Value = 44453
$ws.Range("A15").NumberFormat = "d-mmm"
$ws.Range("B15").Value = "4.5 hours"
$ws.Range("C15").Value = "Task: Updated the database procedures with try catch blocks and added transactions there, also added a new table that has a dataState value in it for updates etc, implemented the refresh on the controller side and hooked an event to it so that the UI can listen to that, removed transactions from the db class"

$ws.Range("A16").Select()
